# Update column G (header "K", formerly populated from "Strike#") with
# regenerated strikeout counts for each game row, per the source data
# regeneration described in the commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 2
    3  = 5
    4  = 11
    5  = 5
    6  = 5
    7  = 2
    8  = 3
    9  = 6
    10 = 1
    11 = 8
    12 = 1
    13 = 5
    14 = 4
    15 = 4
    16 = 7
    17 = 8
    18 = 3
    19 = 7
    20 = 5
    21 = 6
    22 = 7
    23 = 1
    24 = 6
    25 = 2
    26 = 3
    27 = 5
    28 = 3
    29 = 5
    30 = 2
    31 = 4
    32 = 4
    33 = 2
}

foreach ($row in $newValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $newValues[$row]
}
